$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Bump the ticket/date serial in A1 by one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Update the price in D29 (1570 -> 960)
$ws.Range("D29").Value = 960
